$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing sheet's header row
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy date-format styling (column A) from an existing sheet
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats

# Data rows
$newSheet.Range("A2").Value = 45564.99999999999
$newSheet.Range("B2").Value = 406
$newSheet.Range("C2").Value = 232.9555113933152
$newSheet.Range("D2").Value = 583.8843976535485

$newSheet.Range("A3").Value = 45578.99999999999
$newSheet.Range("B3").Value = 255
$newSheet.Range("C3").Value = 75.34611325312835
$newSheet.Range("D3").Value = 429.7924488023465

$newSheet.Range("A4").Value = 45592.99999999999
$newSheet.Range("B4").Value = 105
$newSheet.Range("C4").Value = -69.69374535678364
$newSheet.Range("D4").Value = 278.1072579897799

$newSheet.Range("A5").Value = 45599.99999999999
$newSheet.Range("B5").Value = 30
$newSheet.Range("C5").Value = -146.1070819346028
$newSheet.Range("D5").Value = 221.1980598631854

$newSheet.Range("A6").Value = 45606.99999999999
$newSheet.Range("B6").Value = 0
$newSheet.Range("C6").Value = -224.1524298826996
$newSheet.Range("D6").Value = 130.1052568253197

$newSheet.Range("A7").Value = 45613.99999999999
$newSheet.Range("B7").Value = 0
$newSheet.Range("C7").Value = -292.3763898003348
$newSheet.Range("D7").Value = 54.84718365461671

$newSheet.Range("A8").Value = 45620.99999999999
$newSheet.Range("B8").Value = 0
$newSheet.Range("C8").Value = -387.6210803131578
$newSheet.Range("D8").Value = -13.91588591418187

$newSheet.Range("A9").Value = 45627.99999999999
$newSheet.Range("B9").Value = 0
$newSheet.Range("C9").Value = -456.3717771817866
$newSheet.Range("D9").Value = -97.31561967974737

$newSheet.Range("A10").Value = 45634.99999999999
$newSheet.Range("B10").Value = 0
$newSheet.Range("C10").Value = -510.4310483287725
$newSheet.Range("D10").Value = -161.4629452169551

$newSheet.Range("A11").Value = 45641.99999999999
$newSheet.Range("B11").Value = 0
$newSheet.Range("C11").Value = -589.3150527193153
$newSheet.Range("D11").Value = -237.8802165356711

$newSheet.Range("A12").Value = 45648.99999999999
$newSheet.Range("B12").Value = 0
$newSheet.Range("C12").Value = -673.9774858750671
$newSheet.Range("D12").Value = -306.9087785964843

$newSheet.Range("A13").Value = 45655.99999999999
$newSheet.Range("B13").Value = 0
$newSheet.Range("C13").Value = -753.9998917165079
$newSheet.Range("D13").Value = -390.4919645754753

$newSheet.Range("A14").Value = 45662.99999999999
$newSheet.Range("B14").Value = 0
$newSheet.Range("C14").Value = -827.8417318064276
$newSheet.Range("D14").Value = -443.1008968363563

$newSheet.Range("A1").Select()
